$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 669, pushing the existing
# rows 669-762 down to 671-764 (dimension grows from R762 to R764).
$ws.Range("669:670").EntireRow.Insert()

# Give the new Fecha (date) cells the same number format as the
# surrounding date column (style index 2 in the original workbook).
$ws.Range("D669").NumberFormat = $ws.Range("D671").NumberFormat
$ws.Range("D670").NumberFormat = $ws.Range("D671").NumberFormat

# --- New row 669 (Primera) ---
$ws.Range("A669").Value = 3
$ws.Range("B669").Value = "Femacal de La Calera"
$ws.Range("C669").Value = "Coquimbo"
$ws.Range("D669").Value = 44776
$ws.Range("E669").Value = 5
$ws.Range("F669").Value = 100114014
$ws.Range("G669").Value = "Betarraga"
$ws.Range("H669").Value = "Sin especificar"
$ws.Range("I669").Value = "Primera"
$ws.Range("J669").Value = 3700
$ws.Range("K669").Value = 850
$ws.Range("L669").Value = 900
$ws.Range("M669").Value = 874
$ws.Range("N669").Value = "`$/paquete 4 unidades"
$ws.Range("O669").Value = "Provincia de Quillota"
$ws.Range("P669").Value = 218
$ws.Range("Q669").Value = 4
$ws.Range("R669").Value = "Hortaliza"

# --- New row 670 (Segunda) ---
$ws.Range("A670").Value = 3
$ws.Range("B670").Value = "Femacal de La Calera"
$ws.Range("C670").Value = "Coquimbo"
$ws.Range("D670").Value = 44776
$ws.Range("E670").Value = 5
$ws.Range("F670").Value = 100114014
$ws.Range("G670").Value = "Betarraga"
$ws.Range("H670").Value = "Sin especificar"
$ws.Range("I670").Value = "Segunda"
$ws.Range("J670").Value = 1500
$ws.Range("K670").Value = 700
$ws.Range("L670").Value = 700
$ws.Range("M670").Value = 700
$ws.Range("N670").Value = "`$/paquete 4 unidades"
$ws.Range("O670").Value = "Provincia de Quillota"
$ws.Range("P670").Value = 175
$ws.Range("Q670").Value = 4
$ws.Range("R670").Value = "Hortaliza"
